$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

# The "Individuals" sheet has header rows 1-2 and data rows 3-29.
# Column H = phone_number_1, Column I = phone_number_2.
# Replace the placeholder phone numbers with realistic ones, alternating
# by row parity (odd data rows vs even data rows), matching rows 3/4.

For ($row = 3; $row -le 29; $row++) {
    if (($row % 2) -eq 1) {
        $ws.Cells.Item($row, 8).Value = "+44 1632 960852"
        $ws.Cells.Item($row, 9).Value = "+1-541-754-3010"
    } else {
        $ws.Cells.Item($row, 8).Value = "+1-613-555-0182"
        $ws.Cells.Item($row, 9).Value = "+36 55 979 922"
    }
}
